$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new rows of data (rows 6 and 7), copying the formatting from
#     the last existing data row (row 5) so the new cells pick up style s=2 ---
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 6: id=5, nome="sd", valor=3, data=2025-06-01
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "sd"
$ws.Range("C6").Value = 3

# Row 7: id=6, nome="sd", valor=2, data=2025-06-01
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "sd"
$ws.Range("C7").Value = 2

# Copy the "data" text value from an existing cell so it is stored as text
# (matching the source data) instead of being auto-converted to a date.
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Update the saved selection to match the target view state
$ws.Range("D9").Select()

$wb.Save()
